# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker table (B15:J15 header, rows 16-36) previously listed each of the
# 3 trabajadores with their 7 "Periodo Mora" rows grouped by person
# (2107 down to 2101). This edit re-sorts the same 21 data rows so they are
# grouped by period instead (ascending 2101 -> 2107), with the 3 workers
# appearing in the same relative order inside every period group. Pedro's
# "Salario Basico" is also brought in line with the other two workers
# (738000 -> 877803) now that the periods are interleaved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ Row = 16; Doc = "CC"; Num = "73350956"; Nombre = "PEDRO ANTONIO PATERNINA CORONADO"; Periodo = "2101"; Mora = 35112; Salario = 877803 },
    @{ Row = 17; Doc = "CC"; Num = "9076063";  Nombre = "ORLANDO CABARCAS GUZMAN";          Periodo = "2101"; Mora = 35112; Salario = 877803 },
    @{ Row = 18; Doc = "CC"; Num = "73189894"; Nombre = "WINDER PACHECO RIVERO";            Periodo = "2101"; Mora = 35112; Salario = 877803 },

    @{ Row = 19; Doc = "CC"; Num = "73350956"; Nombre = "PEDRO ANTONIO PATERNINA CORONADO"; Periodo = "2102"; Mora = 35112; Salario = 877803 },
    @{ Row = 20; Doc = "CC"; Num = "9076063";  Nombre = "ORLANDO CABARCAS GUZMAN";          Periodo = "2102"; Mora = 35112; Salario = 877803 },
    @{ Row = 21; Doc = "CC"; Num = "73189894"; Nombre = "WINDER PACHECO RIVERO";            Periodo = "2102"; Mora = 35112; Salario = 877803 },

    @{ Row = 22; Doc = "CC"; Num = "73350956"; Nombre = "PEDRO ANTONIO PATERNINA CORONADO"; Periodo = "2103"; Mora = 35112; Salario = 877803 },
    @{ Row = 23; Doc = "CC"; Num = "9076063";  Nombre = "ORLANDO CABARCAS GUZMAN";          Periodo = "2103"; Mora = 35112; Salario = 877803 },
    @{ Row = 24; Doc = "CC"; Num = "73189894"; Nombre = "WINDER PACHECO RIVERO";            Periodo = "2103"; Mora = 35112; Salario = 877803 },

    @{ Row = 25; Doc = "CC"; Num = "73350956"; Nombre = "PEDRO ANTONIO PATERNINA CORONADO"; Periodo = "2104"; Mora = 35112; Salario = 877803 },
    @{ Row = 26; Doc = "CC"; Num = "9076063";  Nombre = "ORLANDO CABARCAS GUZMAN";          Periodo = "2104"; Mora = 35112; Salario = 877803 },
    @{ Row = 27; Doc = "CC"; Num = "73189894"; Nombre = "WINDER PACHECO RIVERO";            Periodo = "2104"; Mora = 35112; Salario = 877803 },

    @{ Row = 28; Doc = "CC"; Num = "73350956"; Nombre = "PEDRO ANTONIO PATERNINA CORONADO"; Periodo = "2105"; Mora = 35112; Salario = 877803 },
    @{ Row = 29; Doc = "CC"; Num = "9076063";  Nombre = "ORLANDO CABARCAS GUZMAN";          Periodo = "2105"; Mora = 35112; Salario = 877803 },
    @{ Row = 30; Doc = "CC"; Num = "73189894"; Nombre = "WINDER PACHECO RIVERO";            Periodo = "2105"; Mora = 35112; Salario = 877803 },

    @{ Row = 31; Doc = "CC"; Num = "73350956"; Nombre = "PEDRO ANTONIO PATERNINA CORONADO"; Periodo = "2106"; Mora = 35112; Salario = 877803 },
    @{ Row = 32; Doc = "CC"; Num = "9076063";  Nombre = "ORLANDO CABARCAS GUZMAN";          Periodo = "2106"; Mora = 35112; Salario = 877803 },
    @{ Row = 33; Doc = "CC"; Num = "73189894"; Nombre = "WINDER PACHECO RIVERO";            Periodo = "2106"; Mora = 35112; Salario = 877803 },

    @{ Row = 34; Doc = "CC"; Num = "73350956"; Nombre = "PEDRO ANTONIO PATERNINA CORONADO"; Periodo = "2107"; Mora = 29260; Salario = 877803 },
    @{ Row = 35; Doc = "CC"; Num = "9076063";  Nombre = "ORLANDO CABARCAS GUZMAN";          Periodo = "2107"; Mora = 29260; Salario = 877803 },
    @{ Row = 36; Doc = "CC"; Num = "73189894"; Nombre = "WINDER PACHECO RIVERO";            Periodo = "2107"; Mora = 29260; Salario = 877803 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("B$row").Value = $r.Doc
    $ws.Range("C$row").Value = $r.Num
    $ws.Range("D$row").Value = $r.Nombre
    $ws.Range("E$row").Value = $r.Periodo
    $ws.Range("F$row").Value = $r.Mora
    $ws.Range("G$row").Value = $r.Salario
}
